$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1095
$ws.Range("F6").Value = 3338
$ws.Range("F11").Value = 578
$ws.Range("F16").Value = 1708
$ws.Range("F17").Value = 1708
$ws.Range("F19").Value = 337
$ws.Range("F25").Value = 668
$ws.Range("F26").Value = 77806
$ws.Range("F27").Value = 77806
$ws.Range("F32").Value = 475
$ws.Range("F37").Value = 943
$ws.Range("F38").Value = 270
$ws.Range("F41").Value = 64
$ws.Range("F42").Value = 1169
$ws.Range("F43").Value = 5419

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F23").Value = 498
$ws.Range("F24").Value = 498
$ws.Range("F43").Value = 19

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 707
$ws.Range("F6").Value = 567

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 707
$ws.Range("F6").Value = 1095
$ws.Range("F10").Value = 3338
$ws.Range("F18").Value = 567
$ws.Range("F20").Value = 578
$ws.Range("F23").Value = 1708
$ws.Range("F24").Value = 1708
$ws.Range("F26").Value = 337
$ws.Range("F32").Value = 668
$ws.Range("F33").Value = 77806
$ws.Range("F36").Value = 475
$ws.Range("F40").Value = 498
$ws.Range("F43").Value = 270
$ws.Range("F47").Value = 5419
